$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell -> new text value, derived from the authoritative diff.
$updates = @{
    'D2' = '62.227.14'
    'E2' = '  +3.02%  '
    'D3' = '3.406.61'
    'E3' = '  +3.34%  '
    'E4' = '  +0.07%  '
    'D5' = '406.97'
    'E5' = '  -0.37%  '
    'D6' = '130.98'
    'E6' = '  +16.45%  '
    'E7' = '  +6.75%  '
    'E8' = '  +0.05%  '
    'E9' = '  +9.09%  '
    'E10' = '  +11.08%  '
    'D11' = '42.25'
    'E11' = '  +8.40%  '
    'E12' = '  -0.38%  '
    'D13' = '3.960.59'
    'E13' = '  +3.88%  '
    'D14' = '8.56'
    'E14' = '  +4.65%  '
    'E15' = '  +4.02%  '
    'D16' = '3.414.16'
    'E16' = '  +3.66%  '
    'B17' = 'WrappedBTC'
    'C17' = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
    'D17' = '62.077.60'
    'E17' = '  +3.23%  '
    'B18' = 'Uniswap'
    'C18' = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
    'D18' = '11.54'
    'E18' = '  +7.06%  '
    'E19' = '  +4.62%  '
    'E20' = '  +17.43%  '
    'E21' = '  -0.80%  '
    'D22' = '82.86'
    'E22' = '  +12.61%  '
    'D23' = '13.20'
    'E23' = '  +6.28%  '
    'D24' = '307.99'
    'E24' = '  +4.05%  '
    'E25' = '  +2.63%  '
    'D26' = '8.59'
    'E26' = '  +14.80%  '
    'B27' = 'LEO'
    'C27' = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
    'D27' = '4.69'
    'E27' = '  +9.74%  '
    'B28' = 'EthereumClassic'
    'C28' = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
    'D28' = '29.80'
    'E28' = '  +2.16%  '
    'B29' = 'RenderToken'
    'C29' = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
    'D29' = '7.49'
    'E29' = '  +1.76%  '
    'B30' = 'Kaspa'
    'C30' = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
    'D30' = '0.175'
    'E30' = '  +1.55%  '
    'E31' = '  +2.82%  '
    'D32' = '11.79'
    'E32' = '  +5.37%  '
    'D33' = '2.63'
    'E33' = '  +5.99%  '
    'D34' = '42.53'
    'E34' = '  +8.07%  '
    'D35' = '0.999'
    'E35' = '  -0.11%  '
    'E36' = '  +1.08%  '
    'E37' = '  +0.80%  '
    'D38' = '0.997'
    'E38' = '  -0.04%  '
    'E39' = '  +3.57%  '
    'D40' = '2.99'
    'E40' = '  -2.75%  '
    'E41' = '  +8.41%  '
    'E42' = '  +4.91%  '
    'D43' = '137.38'
    'E43' = '  +2.78%  '
    'D44' = '3.98'
    'E44' = '  +5.10%  '
    'D45' = '17.09'
    'E45' = '  +5.05%  '
    'E46' = '  -3.55%  '
    'E47' = '  +2.26%  '
    'D48' = '21.85'
    'E48' = '  +4.56%  '
    'D49' = '2.155.29'
    'E49' = '  +1.41%  '
    'D50' = '3.746.18'
    'E50' = '  +3.80%  '
    'E51' = '  -0.58%  '
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    # Force text format so numeric-looking strings (e.g. "13.20", "0.175",
    # "62.227.14") are stored verbatim as text instead of being coerced into
    # floating point numbers (which would drop formatting like trailing
    # zeros or introduce binary rounding noise).
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
}
